$d = $word.ActiveDocument

# 1. Resize the table grid columns (792/554/6573 twips -> 800/560/6560 twips).
#    Word's Column.Width is expressed in points (1 pt = 20 twips).
$t = $d.Tables.Item(1)
$t.Columns.Item(1).Width = 800 / 20
$t.Columns.Item(2).Width = 560 / 20
$t.Columns.Item(3).Width = 6560 / 20

# 2. Update the solution text for problem 4 (replace the discrete-random-variable
#    explanation with the probability-of-1/probability-of-0 explanation).
$d.Content.Find.Execute(
    "A discrete random variable is something that varies following a specific pattern",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "A probability of 1 implies an event is certain to happen. A probability of 0",
    2)

$d.Content.Find.Execute(
    "or distribution over the long run. They are discrete if they can be listed.",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "implies it is impossible to happen, or certain to not happen.",
    2)
